$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format on columns that receive numeric-looking
# string values (D, E, G) so Excel does not silently convert them to real
# numbers/percentages. ClearFormats() afterwards removes the temporary style
# so the resulting cells have no style attribute, matching the original layout.
# (kept as two single-area ranges; a multi-area "D2:E51,G2:G51" union range
# was observed to let ClearFormats silently revert G back to a real number)
$rngDE = $ws.Range("D2:E51")
$rngG = $ws.Range("G2:G51")
$rngDE.NumberFormat = "@"
$rngG.NumberFormat = "@"

$ws.Range("D2").Value = '300.92'
$ws.Range("E2").Value = '2.38%'
$ws.Range("G2").Value = '12'
$ws.Range("D3").Value = '42.57'
$ws.Range("E3").Value = '5.03%'
$ws.Range("G3").Value = '12'
$ws.Range("D4").Value = '4.991'
$ws.Range("E4").Value = '-0.45%'
$ws.Range("G4").Value = '12'
$ws.Range("D5").Value = '0.07657'
$ws.Range("E5").Value = '3.34%'
$ws.Range("G5").Value = '12'
$ws.Range("D6").Value = '1.607'
$ws.Range("E6").Value = '2.02%'
$ws.Range("G6").Value = '12'
$ws.Range("D7").Value = '0.9922'
$ws.Range("E7").Value = '7.21%'
$ws.Range("G7").Value = '12'
$ws.Range("E8").Value = '1.46%'
$ws.Range("G8").Value = '12'
$ws.Range("D9").Value = '0.1206'
$ws.Range("E9").Value = '0.03%'
$ws.Range("G9").Value = '12'
$ws.Range("D10").Value = '0.1842'
$ws.Range("E10").Value = '1.52%'
$ws.Range("G10").Value = '12'
$ws.Range("D11").Value = '0.09019'
$ws.Range("E11").Value = '2.88%'
$ws.Range("G11").Value = '12'
$ws.Range("D12").Value = '0.04112'
$ws.Range("E12").Value = '-6.18%'
$ws.Range("G12").Value = '12'
$ws.Range("D13").Value = '0.1047'
$ws.Range("E13").Value = '-0.83%'
$ws.Range("G13").Value = '12'
$ws.Range("D14").Value = '0.001279'
$ws.Range("E14").Value = '0.02%'
$ws.Range("G14").Value = '12'
$ws.Range("D15").Value = '0.005863'
$ws.Range("E15").Value = '-2.74%'
$ws.Range("G15").Value = '12'
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = '0.007435'
$ws.Range("E16").Value = '1,895.55%'
$ws.Range("G16").Value = '12'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.348'
$ws.Range("E17").Value = '0.19%'
$ws.Range("G17").Value = '12'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '4.394'
$ws.Range("E18").Value = '2.31%'
$ws.Range("G18").Value = '12'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3335'
$ws.Range("E19").Value = '0.54%'
$ws.Range("G19").Value = '12'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").Value = '8.352'
$ws.Range("E20").Value = '6.36%'
$ws.Range("G20").Value = '12'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '0.1379'
$ws.Range("E21").Value = '-0.83%'
$ws.Range("G21").Value = '12'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").Value = '0.3284'
$ws.Range("E22").Value = '13.89%'
$ws.Range("G22").Value = '12'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").Value = '0.04125'
$ws.Range("E23").Value = '4.89%'
$ws.Range("G23").Value = '12'
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").Value = '0.001261'
$ws.Range("E24").Value = '0.12%'
$ws.Range("G24").Value = '12'
$ws.Range("B25").Value = 'HotbitToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D25").Value = '0.003958'
$ws.Range("E25").Value = '4.60%'
$ws.Range("G25").Value = '12'
$ws.Range("B26").Value = 'NitroEx'
$ws.Range("C26").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D26").Value = '0.0001343'
$ws.Range("E26").Value = '9.25%'
$ws.Range("G26").Value = '12'
$ws.Range("B27").Value = 'Spectre.aiUtilityToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range("G27").Value = '12'
$ws.Range("B28").Value = 'LegolasExchange'
$ws.Range("C28").Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range("G28").Value = '12'
$ws.Range("B29").Value = 'BitZToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range("G29").Value = '12'
$ws.Range("B30").Value = 'Birake'
$ws.Range("C30").Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range("G30").Value = '12'
$ws.Range("B31").Value = 'NashExchange'
$ws.Range("C31").Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
$ws.Range("G31").Value = '12'
$ws.Range("B32").Value = 'AAXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("G32").Value = '12'
$ws.Range("B33").Value = 'CenX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
$ws.Range("G33").Value = '12'
$ws.Range("B34").Value = 'BNIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
$ws.Range("G34").Value = '12'
$ws.Range("G35").Value = '12'
$ws.Range("G36").Value = '12'
$ws.Range("G37").Value = '12'
$ws.Range("D38").Value = '0.02447'
$ws.Range("E38").Value = '5.14%'
$ws.Range("G38").Value = '12'
$ws.Range("D39").Value = '0.05263'
$ws.Range("E39").Value = '3.40%'
$ws.Range("G39").Value = '12'
$ws.Range("D40").Value = '0.006420'
$ws.Range("E40").Value = '1.61%'
$ws.Range("G40").Value = '12'
$ws.Range("D41").Value = '0.007641'
$ws.Range("E41").Value = '-2.11%'
$ws.Range("G41").Value = '12'
$ws.Range("D42").Value = '0.1343'
$ws.Range("E42").Value = '4.08%'
$ws.Range("G42").Value = '12'
$ws.Range("D43").Value = '0.007337'
$ws.Range("E43").Value = '-0.75%'
$ws.Range("G43").Value = '12'
$ws.Range("D44").Value = '0.007294'
$ws.Range("E44").Value = '0.11%'
$ws.Range("G44").Value = '12'
$ws.Range("D45").Value = '0.3018'
$ws.Range("E45").Value = '3.38%'
$ws.Range("G45").Value = '12'
$ws.Range("D46").Value = '0.00006427'
$ws.Range("E46").Value = '3.43%'
$ws.Range("G46").Value = '12'
$ws.Range("D47").Value = '0.00000000744'
$ws.Range("E47").Value = '-0.74%'
$ws.Range("G47").Value = '12'
$ws.Range("D48").Value = '0.04596'
$ws.Range("E48").Value = '-1.84%'
$ws.Range("G48").Value = '12'
$ws.Range("E49").Value = '-0.08%'
$ws.Range("G49").Value = '12'
$ws.Range("D50").Value = '0.00002084'
$ws.Range("E50").Value = '-0.74%'
$ws.Range("G50").Value = '12'
$ws.Range("D51").Value = '0.0001985'
$ws.Range("E51").Value = '-0.74%'
$ws.Range("G51").Value = '12'

$rngDE.ClearFormats()
$rngG.ClearFormats()

